$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Majors")

# Rename "All Majors Combined" to "Bachelor's in All Majors Combined"
$ws.Range("A2").Value = "Bachelor's in All Majors Combined"

# Reflect the active cell left selected after the edit (A3)
$ws.Range("A3").Select()
